# Remove outliers with extreme offsets
#
# The underlying author edit (per the target diff) inserts a new data row
# into the "Tabelle1" worksheet just above the "FreeHead" / "Reference"
# block (i.e. at sheet row 19), pushing the existing rows 19-28 down to
# 20-29. The new row documents a "General" / "accuracy" measurement with
# an offset of 0.5 (d=50) and a text note of "e" in the d=50cm-offset
# column — this is the note referenced by "Remove outliers with extreme
# offsets" (flagging the row as excluded/annotated rather than deleting
# historical data wholesale).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Insert a new blank row at row 19; Excel pushes the old row 19 ("FreeHead",
# C=1, ...) and everything below it down by one row, and carries the
# formatting of the row above (row 18, the "General"/"accuracy" row) into
# the freshly inserted row 19 -- which matches the target styles (s="4" for
# A/B/D, s="12" for C/E/F, s="19" for G).
$ws.Rows.Item(19).Insert()

# The row-insert default formatting also spills into columns J and K
# (copied from row 18's J/K cells); the new row should not carry any
# content/formatting there, so clear those two cells back out.
$ws.Range("J19:K19").Clear()

# Populate the new row's contents.
$ws.Cells.Item(19, 1).Value = "General"
$ws.Cells.Item(19, 2).Value = "accuracy"
$ws.Cells.Item(19, 3).Value = 0.5
$ws.Cells.Item(19, 4).Value = 50
$ws.Cells.Item(19, 6).Value = "e"

# Selection / view bookkeeping (best effort; mirrors the author ending up
# with F19 selected after the edit).
$ws.Range("F19").Select()
